$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows (2-7) hold weekly Pomelo price records. The data is being
# reordered (rows shuffled) while columns A,B,C,E,F,G,H,I,J,K,L,R,T stay
# constant across all of them. Only D (Fecha), M (Volumen), N (Precio
# minimo), O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de
# comercializacion) and S (Precio $/Kg) differ per row.

$rows = @(
    @{ Row = 2; D = 44400; M = 100; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos";            S = 714 },
    @{ Row = 3; D = 44397; M = 60;  N = 11000; O = 11000; P = 11000; Q = "`$/caja 14 kilos";            S = 786 },
    @{ Row = 4; D = 44351; M = 300; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada";  S = 714 },
    @{ Row = 5; D = 44309; M = 300; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada";  S = 500 },
    @{ Row = 6; D = 44162; M = 120; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada";  S = 500 },
    @{ Row = 7; D = 44176; M = 250; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada";  S = 500 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $r.N   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $r.O   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $r.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($row, 19).Value = $r.S   # S: Precio $/Kg
}
